$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot (includes a coin-order swap at rows 50-51)

$ws.Range("D2").Value = '26.659.67'
$ws.Range("E2").Value = '  +1.04%  '

$ws.Range("D3").Value = '1.633.17'
$ws.Range("E3").Value = '  +1.09%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.489'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.62%  '

$ws.Range("E8").Value = '  +0.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0835'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.62%  '

$ws.Range("D12").Value = '1.864.13'
$ws.Range("E12").Value = '  +1.32%  '

$ws.Range("D13").Value = '1.630.32'
$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("E14").Value = '  +0.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '

$ws.Range("D16").Value = '26.660.10'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.43%  '

$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.99'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.98%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '

$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("E27").Value = '  -1.11%  '

$ws.Range("E28").Value = '  +0.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0516'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.92%  '

$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("E34").Value = '  +0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("D36").Value = '1.163.53'
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("E37").Value = '  +1.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.806'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("E39").Value = '  +0.04%  '

$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("E42").Value = '  +1.00%  '

$ws.Range("E43").Value = '  +1.99%  '

$ws.Range("D44").Value = '1.773.08'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.76%  '

$ws.Range("E46").Value = '  +0.95%  '

$ws.Range("E47").Value = '  +7.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.34%  '

$ws.Range("E49").Value = '  +0.72%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.410'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.80%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.14%  '
